# "Membres du groupe" roster was uploaded into the previously-empty Sheet1:
# a two-column (Prenoms / Nom) table of five rows (header + four members).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Prénoms"
$ws.Range("B1").Value = "Nom"

$ws.Range("A2").Value = "Khadidiatou"
$ws.Range("B2").Value = "Coulibaly"

$ws.Range("A3").Value = "Tamsir"
$ws.Range("B3").Value = "Ndong"

$ws.Range("A4").Value = "Samba"
$ws.Range("B4").Value = "Dieng"

$ws.Range("A5").Value = "Jeanne De La Flèche"
$ws.Range("B5").Value = "Onanena Amana"

# Column A was widened (best-fit) to accommodate the longest first name.
$ws.Columns.Item(1).AutoFit() | Out-Null

# The author's last selection before saving was cell M9.
$ws.Range("M9").Select() | Out-Null
